# feat: add 2022-Q3 data
#
# Before: sheets = [总计, 2022-Q2]
# After:  sheets = [总计, 2022-Q3, 2022-Q2]
#   - A new "2022-Q3" sheet is inserted right after "总计", taking over the
#     worksheet part that used to be "2022-Q2" (same position / relationship id).
#   - The former "2022-Q2" data is preserved verbatim in a brand-new sheet that
#     is placed after "2022-Q3" and renamed back to "2022-Q2".
#   - The "总计" summary sheet gets a new row for 2022-Q3 (pushing the old
#     2022-Q2 row down to row 3).

$wb = $excel.ActiveWorkbook
$wsTotal = $wb.Worksheets.Item("总计")
$wsQ3 = $wb.Worksheets.Item("2022-Q2")

# ------------------------------------------------------------------
# 1) Preserve the existing "2022-Q2" worksheet's formatting by copying
#    it into a fresh worksheet before we overwrite anything, then rename
#    the original sheet to "2022-Q3" (so it keeps its original position /
#    relationship id, matching how the new quarter sheet should be wired).
# ------------------------------------------------------------------
$wsQ3.Range("B1:H1").Copy()

$wsQ3.Name = "2022-Q3"

$wsQ2 = $wb.Worksheets.Add($null, $wsQ3)
$wsQ2.Name = "2022-Q2"
$wsQ2.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsQ2.Range("B1").Value = "基金代码"
$wsQ2.Range("C1").Value = "基金名称"
$wsQ2.Range("D1").Value = "基金规模"
$wsQ2.Range("E1").Value = "股票总仓位"
$wsQ2.Range("F1").Value = "仓位占比"
$wsQ2.Range("G1").Value = "持有市值(亿元)"
$wsQ2.Range("H1").Value = "仓位排名"

$wsQ2.Range("B1").Copy()
$wsQ2.Range("A2:A3").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$wsQ2.Range("A2").Value = 0
$wsQ2.Range("A3").Value = 1

$oldTxt = $wsQ2.Range("B2:G3")
$oldTxt.NumberFormat = "@"
$wsQ2.Range("B2").Value = "005493"
$wsQ2.Range("C2").Value = "鑫元价值精选灵活配置混合A"
$wsQ2.Range("D2").Value = "0.63"
$wsQ2.Range("E2").Value = "60.60"
$wsQ2.Range("F2").Value = "2.45"
$wsQ2.Range("G2").Value = "0.0154"
$wsQ2.Range("B3").Value = "005494"
$wsQ2.Range("C3").Value = "鑫元价值精选灵活配置混合C"
$wsQ2.Range("D3").Value = "0.00"
$wsQ2.Range("E3").Value = "60.60"
$wsQ2.Range("F3").Value = "2.45"
$oldTxt.ClearFormats()

$wsQ2.Range("G3").Value = 0
$wsQ2.Range("H2").Value = 8
$wsQ2.Range("H3").Value = 8

# ------------------------------------------------------------------
# 2) Rebuild the "2022-Q3" worksheet (previously "2022-Q2") with the
#    new quarter's fund holdings.
# ------------------------------------------------------------------
$wsQ3.Cells.Clear()

$wsTotal.Range("B1").Copy()
$wsQ3.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsQ3.Range("B1").Value = "基金代码"
$wsQ3.Range("C1").Value = "基金名称"
$wsQ3.Range("D1").Value = "基金规模"
$wsQ3.Range("E1").Value = "股票总仓位"
$wsQ3.Range("F1").Value = "仓位占比"
$wsQ3.Range("G1").Value = "持有市值(亿元)"
$wsQ3.Range("H1").Value = "仓位排名"

$wsQ3.Range("B1").Copy()
$wsQ3.Range("A2:A9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$q3Data = @(
  @("011815", "恒越优势精选混合",         "3.22", "90.91", "3.89", "0.1253", 2),
  @("013721", "信澳景气优选混合A",        "1.06", "82.64", "4.17", "0.0442", 10),
  @("013028", "恒越品质生活混合",         "1.35", "89.92", "3.19", "0.0431", 6),
  @("011629", "银河核心优势混合",         "2.68", "26.95", "1.00", "0.0268", 10),
  @("013722", "信澳景气优选混合C",        "0.47", "82.64", "4.17", "0.0196", 10),
  @("710002", "富安达策略精选混合",       "0.59", "50.67", "1.74", "0.0103", 6),
  @("002584", "富安达长盈灵活配置混合A",  "0.11", "46.55", "2.75", "0.0030", 1),
  @("016214", "富安达长盈灵活配置混合C",  "0.02", "46.55", "2.75", "0.0006", 1)
)

$r = 2
foreach ($row in $q3Data) {
  $wsQ3.Range("A$r").Value = ($r - 2)

  $rowTxt = $wsQ3.Range("B$r" + ":G$r")
  $rowTxt.NumberFormat = "@"
  $wsQ3.Range("B$r").Value = $row[0]
  $wsQ3.Range("C$r").Value = $row[1]
  $wsQ3.Range("D$r").Value = $row[2]
  $wsQ3.Range("E$r").Value = $row[3]
  $wsQ3.Range("F$r").Value = $row[4]
  $wsQ3.Range("G$r").Value = $row[5]
  $rowTxt.ClearFormats()

  $wsQ3.Range("H$r").Value = $row[6]

  $r = $r + 1
}

# ------------------------------------------------------------------
# 3) Update the "总计" summary sheet: insert the 2022-Q3 row, pushing the
#    existing 2022-Q2 row down to row 3.
# ------------------------------------------------------------------
$wsTotal.Range("A2:D2").Copy()
$wsTotal.Range("A3:D3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q2"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.02

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 8
$wsTotal.Range("D2").Value = 0.27

$wsTotal.Activate()
$wsTotal.Range("A1").Select()
